$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44159
$ws.Range("K2").Value = 'Castle Brite'
$ws.Range("N2").Value = 14000
$ws.Range("O2").Value = 15000
$ws.Range("P2").Value = 14500
$ws.Range("Q2").Value = '$/caja 15 kilos'
$ws.Range("R2").Value = 'Región Metropolitana'
$ws.Range("S2").Value = 967
$ws.Range("T2").Value = 15
$ws.Range("L3").Value = 'Primera'
$ws.Range("M3").Value = 100
$ws.Range("N3").Value = 19000
$ws.Range("O3").Value = 20000
$ws.Range("P3").Value = 19500
$ws.Range("S3").Value = 1083
$ws.Range("D4").Value = 44559
$ws.Range("K4").Value = 'Modesto'
$ws.Range("L4").Value = 'Segunda'
$ws.Range("M4").Value = 50
$ws.Range("O4").Value = 18000
$ws.Range("P4").Value = 18000
$ws.Range("Q4").Value = '$/caja 18 kilos'
$ws.Range("S4").Value = 1000
$ws.Range("T4").Value = 18
$ws.Range("D5").Value = 44189
$ws.Range("K5").Value = 'Dina'
$ws.Range("L5").Value = 'Primera'
$ws.Range("M5").Value = 200
$ws.Range("N5").Value = 15000
$ws.Range("O5").Value = 16000
$ws.Range("P5").Value = 15500
$ws.Range("Q5").Value = '$/caja 15 kilos granel'
$ws.Range("S5").Value = 1033
$ws.Range("L6").Value = 'Segunda'
$ws.Range("M6").Value = 100
$ws.Range("N6").Value = 14000
$ws.Range("O6").Value = 14000
$ws.Range("P6").Value = 14000
$ws.Range("S6").Value = 933
$ws.Range("D7").Value = 44545
$ws.Range("K7").Value = 'Castle Brite'
$ws.Range("L7").Value = 'Primera'
$ws.Range("N7").Value = 18000
$ws.Range("O7").Value = 19000
$ws.Range("P7").Value = 18500
$ws.Range("Q7").Value = '$/caja 15 kilos'
$ws.Range("S7").Value = 1233
$ws.Range("D8").Value = 44545
$ws.Range("K8").Value = 'Castle Brite'
$ws.Range("L8").Value = 'Segunda'
$ws.Range("M8").Value = 50
$ws.Range("N8").Value = 17000
$ws.Range("O8").Value = 17000
$ws.Range("P8").Value = 17000
$ws.Range("Q8").Value = '$/caja 15 kilos'
$ws.Range("R8").Value = 'Región de O''Higgins'
$ws.Range("S8").Value = 1133
$ws.Range("T8").Value = 15
$ws.Range("D10").Value = 44187
$ws.Range("K10").Value = 'Dina'
$ws.Range("N10").Value = 15000
$ws.Range("O10").Value = 16000
$ws.Range("P10").Value = 15500
$ws.Range("Q10").Value = '$/caja 18 kilos'
$ws.Range("S10").Value = 861
$ws.Range("T10").Value = 18
